# SECOP_exported_1 edit script
# - Fill the new "Licita" column (T) for rows 2-111 with 0/1 flags
# - Apply wrap-text style + taller row height to row 87 (long Q87 text)
# - Trim the trailing blank rows (113-482) down to a single summary row 112
# - Refresh the AutoFilter range / _FilterDatabase defined name / dimension
# - Update the active selection to T112

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Populate column T ("Licita") for the data rows that are missing it ---
$rowsWithZero = @(3,4,6,7,8,9,10,11,12,13,15,16,19,20,24,25,26,28,29,30,31,33,34,35,37,38,39,40,44,46,47,48,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111)
$rowsWithOne = @(49,90)

foreach ($r in $rowsWithZero) {
    $ws.Cells.Item($r, 20).Value = 0
}
foreach ($r in $rowsWithOne) {
    $ws.Cells.Item($r, 20).Value = 1
}

# --- 2. Row 87: Q87 gets wrap-text formatting, row grows to fit the long text ---
$ws.Cells.Item(87, 17).WrapText = $true
$ws.Rows(87).RowHeight = 283.5

# --- 3. Drop the long run of empty trailing rows (113-482), keep one final row ---
$ws.Rows("113:482").Delete()
$ws.Cells.Item(112, 10).ClearContents()
$ws.Cells.Item(112, 20).Value = 1

# --- 4. Rebuild the AutoFilter over the now-shrunk data range ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:T112").AutoFilter() | Out-Null

# --- 5. Keep the _FilterDatabase defined name in sync with the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "SECOP_exported_1!_FilterDatabase") {
        $n.RefersTo = "=SECOP_exported_1!`$A`$1:`$T`$112"
    }
}

# --- 6. Update the visible selection to match the new last cell ---
$ws.Range("T112").Select() | Out-Null

# --- 7. Best-effort: mirror the author's resized/repositioned workbook window ---
$win = $wb.Windows.Item(1)
$win.Left = 20370
$win.Top = -120
$win.Width = 21840
$win.Height = 13020

Write-Output "edit complete"
